$d = $word.ActiveDocument

$replacements = @(
    @{old = "683×6="; new = "776×2="},
    @{old = "376×3="; new = "941×3="},
    @{old = "513×2="; new = "101×7="},
    @{old = "782×7="; new = "201×5="},
    @{old = "801×3="; new = "307×2="},
    @{old = "250×2="; new = "552×7="},
    @{old = "714×4="; new = "467×2="},
    @{old = "833×6="; new = "614×2="},
    @{old = "872×2="; new = "170×5="},
    @{old = "968×2="; new = "905×3="},
    @{old = "450×5="; new = "672×4="},
    @{old = "811×5="; new = "493×8="},
    @{old = "285×4="; new = "505×8="},
    @{old = "197×5="; new = "433×5="},
    @{old = "405×2="; new = "354×7="},
    @{old = "790×5="; new = "911×9="},
    @{old = "344×7="; new = "778×4="},
    @{old = "506×8="; new = "273×8="},
    @{old = "750×2="; new = "821×5="},
    @{old = "719×3="; new = "418×4="},
    @{old = "140×2="; new = "399×4="},
    @{old = "883×7="; new = "132×2="},
    @{old = "373×4="; new = "474×9="},
    @{old = "281×7="; new = "428×6="},
    @{old = "492×8="; new = "854×4="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
